$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 73,18
$arr[0,0] = 4
$arr[0,1] = "Feria Lagunitas de Puerto Montt"
$arr[0,2] = "Los Lagos"
$arr[0,3] = 44645
$arr[0,4] = 10
$arr[0,5] = 100112052
$arr[0,6] = "Albahaca"
$arr[0,7] = "Sin especificar"
$arr[0,8] = "Primera"
$arr[0,9] = 90
$arr[0,10] = 5000
$arr[0,11] = 7000
$arr[0,12] = 6000
$arr[0,13] = "`$/docena de matas"
$arr[0,14] = "Región Metropolitana"
$arr[0,15] = 1000
$arr[0,16] = 6
$arr[0,17] = "Hortaliza"
$arr[1,0] = 4
$arr[1,1] = "Feria Lagunitas de Puerto Montt"
$arr[1,2] = "Los Lagos"
$arr[1,3] = 44222
$arr[1,4] = 10
$arr[1,5] = 100112052
$arr[1,6] = "Albahaca"
$arr[1,7] = "Sin especificar"
$arr[1,8] = "Primera"
$arr[1,9] = 120
$arr[1,10] = 6000
$arr[1,11] = 6000
$arr[1,12] = 6000
$arr[1,13] = "`$/docena de matas"
$arr[1,14] = "Región Metropolitana"
$arr[1,15] = 1000
$arr[1,16] = 6
$arr[1,17] = "Hortaliza"
$arr[2,0] = 4
$arr[2,1] = "Feria Lagunitas de Puerto Montt"
$arr[2,2] = "Los Lagos"
$arr[2,3] = 44383
$arr[2,4] = 10
$arr[2,5] = 100112052
$arr[2,6] = "Albahaca"
$arr[2,7] = "Sin especificar"
$arr[2,8] = "Primera"
$arr[2,9] = 90
$arr[2,10] = 6000
$arr[2,11] = 6000
$arr[2,12] = 6000
$arr[2,13] = "`$/paquete"
$arr[2,14] = "Región de Arica y Parinacota"
$arr[2,15] = 6000
$arr[2,16] = 1
$arr[2,17] = "Hortaliza"
$arr[3,0] = 4
$arr[3,1] = "Feria Lagunitas de Puerto Montt"
$arr[3,2] = "Los Lagos"
$arr[3,3] = 44295
$arr[3,4] = 10
$arr[3,5] = 100112052
$arr[3,6] = "Albahaca"
$arr[3,7] = "Sin especificar"
$arr[3,8] = "Primera"
$arr[3,9] = 90
$arr[3,10] = 7000
$arr[3,11] = 7000
$arr[3,12] = 7000
$arr[3,13] = "`$/docena de matas"
$arr[3,14] = "Región Metropolitana"
$arr[3,15] = 1167
$arr[3,16] = 6
$arr[3,17] = "Hortaliza"
$arr[4,0] = 4
$arr[4,1] = "Feria Lagunitas de Puerto Montt"
$arr[4,2] = "Los Lagos"
$arr[4,3] = 44278
$arr[4,4] = 10
$arr[4,5] = 100112052
$arr[4,6] = "Albahaca"
$arr[4,7] = "Sin especificar"
$arr[4,8] = "Primera"
$arr[4,9] = 120
$arr[4,10] = 6000
$arr[4,11] = 6000
$arr[4,12] = 6000
$arr[4,13] = "`$/docena de matas"
$arr[4,14] = "Región Metropolitana"
$arr[4,15] = 1000
$arr[4,16] = 6
$arr[4,17] = "Hortaliza"
$arr[5,0] = 4
$arr[5,1] = "Feria Lagunitas de Puerto Montt"
$arr[5,2] = "Los Lagos"
$arr[5,3] = 44435
$arr[5,4] = 10
$arr[5,5] = 100112052
$arr[5,6] = "Albahaca"
$arr[5,7] = "Sin especificar"
$arr[5,8] = "Primera"
$arr[5,9] = 170
$arr[5,10] = 6500
$arr[5,11] = 7000
$arr[5,12] = 6765
$arr[5,13] = "`$/paquete"
$arr[5,14] = "Región de Arica y Parinacota"
$arr[5,15] = 6765
$arr[5,16] = 1
$arr[5,17] = "Hortaliza"
$arr[6,0] = 4
$arr[6,1] = "Feria Lagunitas de Puerto Montt"
$arr[6,2] = "Los Lagos"
$arr[6,3] = 44292
$arr[6,4] = 10
$arr[6,5] = 100112052
$arr[6,6] = "Albahaca"
$arr[6,7] = "Sin especificar"
$arr[6,8] = "Primera"
$arr[6,9] = 120
$arr[6,10] = 7000
$arr[6,11] = 7000
$arr[6,12] = 7000
$arr[6,13] = "`$/docena de matas"
$arr[6,14] = "Región Metropolitana"
$arr[6,15] = 1167
$arr[6,16] = 6
$arr[6,17] = "Hortaliza"
$arr[7,0] = 4
$arr[7,1] = "Feria Lagunitas de Puerto Montt"
$arr[7,2] = "Los Lagos"
$arr[7,3] = 44176
$arr[7,4] = 10
$arr[7,5] = 100112052
$arr[7,6] = "Albahaca"
$arr[7,7] = "Sin especificar"
$arr[7,8] = "Primera"
$arr[7,9] = 80
$arr[7,10] = 6500
$arr[7,11] = 6500
$arr[7,12] = 6500
$arr[7,13] = "`$/docena de matas"
$arr[7,14] = "Región Metropolitana"
$arr[7,15] = 1083
$arr[7,16] = 6
$arr[7,17] = "Hortaliza"
$arr[8,0] = 4
$arr[8,1] = "Feria Lagunitas de Puerto Montt"
$arr[8,2] = "Los Lagos"
$arr[8,3] = 44260
$arr[8,4] = 10
$arr[8,5] = 100112052
$arr[8,6] = "Albahaca"
$arr[8,7] = "Sin especificar"
$arr[8,8] = "Primera"
$arr[8,9] = 140
$arr[8,10] = 7000
$arr[8,11] = 7500
$arr[8,12] = 7214
$arr[8,13] = "`$/docena de matas"
$arr[8,14] = "Región Metropolitana"
$arr[8,15] = 1202
$arr[8,16] = 6
$arr[8,17] = "Hortaliza"
$arr[9,0] = 4
$arr[9,1] = "Feria Lagunitas de Puerto Montt"
$arr[9,2] = "Los Lagos"
$arr[9,3] = 44341
$arr[9,4] = 10
$arr[9,5] = 100112052
$arr[9,6] = "Albahaca"
$arr[9,7] = "Sin especificar"
$arr[9,8] = "Primera"
$arr[9,9] = 120
$arr[9,10] = 5000
$arr[9,11] = 5000
$arr[9,12] = 5000
$arr[9,13] = "`$/docena de matas"
$arr[9,14] = "Región Metropolitana"
$arr[9,15] = 833
$arr[9,16] = 6
$arr[9,17] = "Hortaliza"
$arr[10,0] = 4
$arr[10,1] = "Feria Lagunitas de Puerto Montt"
$arr[10,2] = "Los Lagos"
$arr[10,3] = 44245
$arr[10,4] = 10
$arr[10,5] = 100112052
$arr[10,6] = "Albahaca"
$arr[10,7] = "Sin especificar"
$arr[10,8] = "Primera"
$arr[10,9] = 20
$arr[10,10] = 6000
$arr[10,11] = 6000
$arr[10,12] = 6000
$arr[10,13] = "`$/docena de matas"
$arr[10,14] = "Región Metropolitana"
$arr[10,15] = 1000
$arr[10,16] = 6
$arr[10,17] = "Hortaliza"
$arr[11,0] = 4
$arr[11,1] = "Feria Lagunitas de Puerto Montt"
$arr[11,2] = "Los Lagos"
$arr[11,3] = 44532
$arr[11,4] = 10
$arr[11,5] = 100112052
$arr[11,6] = "Albahaca"
$arr[11,7] = "Sin especificar"
$arr[11,8] = "Primera"
$arr[11,9] = 30
$arr[11,10] = 8000
$arr[11,11] = 8000
$arr[11,12] = 8000
$arr[11,13] = "`$/docena de matas"
$arr[11,14] = "Región Metropolitana"
$arr[11,15] = 1333
$arr[11,16] = 6
$arr[11,17] = "Hortaliza"
$arr[12,0] = 4
$arr[12,1] = "Feria Lagunitas de Puerto Montt"
$arr[12,2] = "Los Lagos"
$arr[12,3] = 44638
$arr[12,4] = 10
$arr[12,5] = 100112052
$arr[12,6] = "Albahaca"
$arr[12,7] = "Sin especificar"
$arr[12,8] = "Primera"
$arr[12,9] = 120
$arr[12,10] = 6000
$arr[12,11] = 6000
$arr[12,12] = 6000
$arr[12,13] = "`$/docena de matas"
$arr[12,14] = "Región Metropolitana"
$arr[12,15] = 1000
$arr[12,16] = 6
$arr[12,17] = "Hortaliza"
$arr[13,0] = 4
$arr[13,1] = "Feria Lagunitas de Puerto Montt"
$arr[13,2] = "Los Lagos"
$arr[13,3] = 44442
$arr[13,4] = 10
$arr[13,5] = 100112052
$arr[13,6] = "Albahaca"
$arr[13,7] = "Sin especificar"
$arr[13,8] = "Primera"
$arr[13,9] = 90
$arr[13,10] = 6000
$arr[13,11] = 6000
$arr[13,12] = 6000
$arr[13,13] = "`$/paquete"
$arr[13,14] = "Región de Arica y Parinacota"
$arr[13,15] = 6000
$arr[13,16] = 1
$arr[13,17] = "Hortaliza"
$arr[14,0] = 4
$arr[14,1] = "Feria Lagunitas de Puerto Montt"
$arr[14,2] = "Los Lagos"
$arr[14,3] = 44540
$arr[14,4] = 10
$arr[14,5] = 100112052
$arr[14,6] = "Albahaca"
$arr[14,7] = "Sin especificar"
$arr[14,8] = "Primera"
$arr[14,9] = 130
$arr[14,10] = 7000
$arr[14,11] = 7000
$arr[14,12] = 7000
$arr[14,13] = "`$/docena de matas"
$arr[14,14] = "Región Metropolitana"
$arr[14,15] = 1167
$arr[14,16] = 6
$arr[14,17] = "Hortaliza"
$arr[15,0] = 4
$arr[15,1] = "Feria Lagunitas de Puerto Montt"
$arr[15,2] = "Los Lagos"
$arr[15,3] = 44246
$arr[15,4] = 10
$arr[15,5] = 100112052
$arr[15,6] = "Albahaca"
$arr[15,7] = "Sin especificar"
$arr[15,8] = "Primera"
$arr[15,9] = 110
$arr[15,10] = 6000
$arr[15,11] = 6000
$arr[15,12] = 6000
$arr[15,13] = "`$/docena de matas"
$arr[15,14] = "Región Metropolitana"
$arr[15,15] = 1000
$arr[15,16] = 6
$arr[15,17] = "Hortaliza"
$arr[16,0] = 4
$arr[16,1] = "Feria Lagunitas de Puerto Montt"
$arr[16,2] = "Los Lagos"
$arr[16,3] = 44323
$arr[16,4] = 10
$arr[16,5] = 100112052
$arr[16,6] = "Albahaca"
$arr[16,7] = "Sin especificar"
$arr[16,8] = "Primera"
$arr[16,9] = 120
$arr[16,10] = 8000
$arr[16,11] = 8000
$arr[16,12] = 8000
$arr[16,13] = "`$/docena"
$arr[16,14] = "Región Metropolitana"
$arr[16,15] = 6667
$arr[16,16] = 1.2
$arr[16,17] = "Hortaliza"
$arr[17,0] = 4
$arr[17,1] = "Feria Lagunitas de Puerto Montt"
$arr[17,2] = "Los Lagos"
$arr[17,3] = 44392
$arr[17,4] = 10
$arr[17,5] = 100112052
$arr[17,6] = "Albahaca"
$arr[17,7] = "Sin especificar"
$arr[17,8] = "Primera"
$arr[17,9] = 50
$arr[17,10] = 7000
$arr[17,11] = 7000
$arr[17,12] = 7000
$arr[17,13] = "`$/paquete"
$arr[17,14] = "Región de Arica y Parinacota"
$arr[17,15] = 7000
$arr[17,16] = 1
$arr[17,17] = "Hortaliza"
$arr[18,0] = 4
$arr[18,1] = "Feria Lagunitas de Puerto Montt"
$arr[18,2] = "Los Lagos"
$arr[18,3] = 44589
$arr[18,4] = 10
$arr[18,5] = 100112052
$arr[18,6] = "Albahaca"
$arr[18,7] = "Sin especificar"
$arr[18,8] = "Primera"
$arr[18,9] = 150
$arr[18,10] = 6000
$arr[18,11] = 6000
$arr[18,12] = 6000
$arr[18,13] = "`$/docena de matas"
$arr[18,14] = "Región Metropolitana"
$arr[18,15] = 1000
$arr[18,16] = 6
$arr[18,17] = "Hortaliza"
$arr[19,0] = 4
$arr[19,1] = "Feria Lagunitas de Puerto Montt"
$arr[19,2] = "Los Lagos"
$arr[19,3] = 44546
$arr[19,4] = 10
$arr[19,5] = 100112052
$arr[19,6] = "Albahaca"
$arr[19,7] = "Sin especificar"
$arr[19,8] = "Primera"
$arr[19,9] = 60
$arr[19,10] = 8000
$arr[19,11] = 8000
$arr[19,12] = 8000
$arr[19,13] = "`$/docena de matas"
$arr[19,14] = "Región Metropolitana"
$arr[19,15] = 1333
$arr[19,16] = 6
$arr[19,17] = "Hortaliza"
$arr[20,0] = 4
$arr[20,1] = "Feria Lagunitas de Puerto Montt"
$arr[20,2] = "Los Lagos"
$arr[20,3] = 44316
$arr[20,4] = 10
$arr[20,5] = 100112052
$arr[20,6] = "Albahaca"
$arr[20,7] = "Sin especificar"
$arr[20,8] = "Primera"
$arr[20,9] = 90
$arr[20,10] = 8000
$arr[20,11] = 8000
$arr[20,12] = 8000
$arr[20,13] = "`$/docena de matas"
$arr[20,14] = "Región Metropolitana"
$arr[20,15] = 1333
$arr[20,16] = 6
$arr[20,17] = "Hortaliza"
$arr[21,0] = 4
$arr[21,1] = "Feria Lagunitas de Puerto Montt"
$arr[21,2] = "Los Lagos"
$arr[21,3] = 44582
$arr[21,4] = 10
$arr[21,5] = 100112052
$arr[21,6] = "Albahaca"
$arr[21,7] = "Sin especificar"
$arr[21,8] = "Primera"
$arr[21,9] = 140
$arr[21,10] = 6000
$arr[21,11] = 6000
$arr[21,12] = 6000
$arr[21,13] = "`$/docena de matas"
$arr[21,14] = "Región Metropolitana"
$arr[21,15] = 1000
$arr[21,16] = 6
$arr[21,17] = "Hortaliza"
$arr[22,0] = 4
$arr[22,1] = "Feria Lagunitas de Puerto Montt"
$arr[22,2] = "Los Lagos"
$arr[22,3] = 44159
$arr[22,4] = 10
$arr[22,5] = 100112052
$arr[22,6] = "Albahaca"
$arr[22,7] = "Sin especificar"
$arr[22,8] = "Primera"
$arr[22,9] = 90
$arr[22,10] = 6500
$arr[22,11] = 6500
$arr[22,12] = 6500
$arr[22,13] = "`$/docena de matas"
$arr[22,14] = "Región Metropolitana"
$arr[22,15] = 1083
$arr[22,16] = 6
$arr[22,17] = "Hortaliza"
$arr[23,0] = 4
$arr[23,1] = "Feria Lagunitas de Puerto Montt"
$arr[23,2] = "Los Lagos"
$arr[23,3] = 44386
$arr[23,4] = 10
$arr[23,5] = 100112052
$arr[23,6] = "Albahaca"
$arr[23,7] = "Sin especificar"
$arr[23,8] = "Primera"
$arr[23,9] = 80
$arr[23,10] = 7000
$arr[23,11] = 7000
$arr[23,12] = 7000
$arr[23,13] = "`$/paquete"
$arr[23,14] = "Región de Arica y Parinacota"
$arr[23,15] = 7000
$arr[23,16] = 1
$arr[23,17] = "Hortaliza"
$arr[24,0] = 4
$arr[24,1] = "Feria Lagunitas de Puerto Montt"
$arr[24,2] = "Los Lagos"
$arr[24,3] = 44320
$arr[24,4] = 10
$arr[24,5] = 100112052
$arr[24,6] = "Albahaca"
$arr[24,7] = "Sin especificar"
$arr[24,8] = "Primera"
$arr[24,9] = 120
$arr[24,10] = 8000
$arr[24,11] = 8000
$arr[24,12] = 8000
$arr[24,13] = "`$/docena de matas"
$arr[24,14] = "Región Metropolitana"
$arr[24,15] = 1333
$arr[24,16] = 6
$arr[24,17] = "Hortaliza"
$arr[25,0] = 4
$arr[25,1] = "Feria Lagunitas de Puerto Montt"
$arr[25,2] = "Los Lagos"
$arr[25,3] = 44204
$arr[25,4] = 10
$arr[25,5] = 100112052
$arr[25,6] = "Albahaca"
$arr[25,7] = "Sin especificar"
$arr[25,8] = "Primera"
$arr[25,9] = 80
$arr[25,10] = 6500
$arr[25,11] = 6500
$arr[25,12] = 6500
$arr[25,13] = "`$/docena de matas"
$arr[25,14] = "Región Metropolitana"
$arr[25,15] = 1083
$arr[25,16] = 6
$arr[25,17] = "Hortaliza"
$arr[26,0] = 4
$arr[26,1] = "Feria Lagunitas de Puerto Montt"
$arr[26,2] = "Los Lagos"
$arr[26,3] = 44313
$arr[26,4] = 10
$arr[26,5] = 100112052
$arr[26,6] = "Albahaca"
$arr[26,7] = "Sin especificar"
$arr[26,8] = "Primera"
$arr[26,9] = 90
$arr[26,10] = 7000
$arr[26,11] = 7000
$arr[26,12] = 7000
$arr[26,13] = "`$/docena de matas"
$arr[26,14] = "Región Metropolitana"
$arr[26,15] = 1167
$arr[26,16] = 6
$arr[26,17] = "Hortaliza"
$arr[27,0] = 4
$arr[27,1] = "Feria Lagunitas de Puerto Montt"
$arr[27,2] = "Los Lagos"
$arr[27,3] = 44460
$arr[27,4] = 10
$arr[27,5] = 100112052
$arr[27,6] = "Albahaca"
$arr[27,7] = "Sin especificar"
$arr[27,8] = "Primera"
$arr[27,9] = 80
$arr[27,10] = 6000
$arr[27,11] = 6000
$arr[27,12] = 6000
$arr[27,13] = "`$/paquete"
$arr[27,14] = "Región de Arica y Parinacota"
$arr[27,15] = 6000
$arr[27,16] = 1
$arr[27,17] = "Hortaliza"
$arr[28,0] = 4
$arr[28,1] = "Feria Lagunitas de Puerto Montt"
$arr[28,2] = "Los Lagos"
$arr[28,3] = 44302
$arr[28,4] = 10
$arr[28,5] = 100112052
$arr[28,6] = "Albahaca"
$arr[28,7] = "Sin especificar"
$arr[28,8] = "Primera"
$arr[28,9] = 80
$arr[28,10] = 8500
$arr[28,11] = 8500
$arr[28,12] = 8500
$arr[28,13] = "`$/docena de matas"
$arr[28,14] = "Región Metropolitana"
$arr[28,15] = 1417
$arr[28,16] = 6
$arr[28,17] = "Hortaliza"
$arr[29,0] = 4
$arr[29,1] = "Feria Lagunitas de Puerto Montt"
$arr[29,2] = "Los Lagos"
$arr[29,3] = 44539
$arr[29,4] = 10
$arr[29,5] = 100112052
$arr[29,6] = "Albahaca"
$arr[29,7] = "Sin especificar"
$arr[29,8] = "Primera"
$arr[29,9] = 60
$arr[29,10] = 8000
$arr[29,11] = 8000
$arr[29,12] = 8000
$arr[29,13] = "`$/docena de matas"
$arr[29,14] = "Región Metropolitana"
$arr[29,15] = 1333
$arr[29,16] = 6
$arr[29,17] = "Hortaliza"
$arr[30,0] = 4
$arr[30,1] = "Feria Lagunitas de Puerto Montt"
$arr[30,2] = "Los Lagos"
$arr[30,3] = 44547
$arr[30,4] = 10
$arr[30,5] = 100112052
$arr[30,6] = "Albahaca"
$arr[30,7] = "Sin especificar"
$arr[30,8] = "Primera"
$arr[30,9] = 90
$arr[30,10] = 8000
$arr[30,11] = 8000
$arr[30,12] = 8000
$arr[30,13] = "`$/docena de matas"
$arr[30,14] = "Región Metropolitana"
$arr[30,15] = 1333
$arr[30,16] = 6
$arr[30,17] = "Hortaliza"
$arr[31,0] = 4
$arr[31,1] = "Feria Lagunitas de Puerto Montt"
$arr[31,2] = "Los Lagos"
$arr[31,3] = 44603
$arr[31,4] = 10
$arr[31,5] = 100112052
$arr[31,6] = "Albahaca"
$arr[31,7] = "Sin especificar"
$arr[31,8] = "Primera"
$arr[31,9] = 140
$arr[31,10] = 6000
$arr[31,11] = 6000
$arr[31,12] = 6000
$arr[31,13] = "`$/docena de matas"
$arr[31,14] = "Región Metropolitana"
$arr[31,15] = 1000
$arr[31,16] = 6
$arr[31,17] = "Hortaliza"
$arr[32,0] = 4
$arr[32,1] = "Feria Lagunitas de Puerto Montt"
$arr[32,2] = "Los Lagos"
$arr[32,3] = 44596
$arr[32,4] = 10
$arr[32,5] = 100112052
$arr[32,6] = "Albahaca"
$arr[32,7] = "Sin especificar"
$arr[32,8] = "Primera"
$arr[32,9] = 120
$arr[32,10] = 6000
$arr[32,11] = 6000
$arr[32,12] = 6000
$arr[32,13] = "`$/docena de matas"
$arr[32,14] = "Región Metropolitana"
$arr[32,15] = 1000
$arr[32,16] = 6
$arr[32,17] = "Hortaliza"
$arr[33,0] = 4
$arr[33,1] = "Feria Lagunitas de Puerto Montt"
$arr[33,2] = "Los Lagos"
$arr[33,3] = 44242
$arr[33,4] = 10
$arr[33,5] = 100112052
$arr[33,6] = "Albahaca"
$arr[33,7] = "Sin especificar"
$arr[33,8] = "Primera"
$arr[33,9] = 30
$arr[33,10] = 6000
$arr[33,11] = 6000
$arr[33,12] = 6000
$arr[33,13] = "`$/docena de matas"
$arr[33,14] = "Región Metropolitana"
$arr[33,15] = 1000
$arr[33,16] = 6
$arr[33,17] = "Hortaliza"
$arr[34,0] = 4
$arr[34,1] = "Feria Lagunitas de Puerto Montt"
$arr[34,2] = "Los Lagos"
$arr[34,3] = 44166
$arr[34,4] = 10
$arr[34,5] = 100112052
$arr[34,6] = "Albahaca"
$arr[34,7] = "Sin especificar"
$arr[34,8] = "Primera"
$arr[34,9] = 90
$arr[34,10] = 6000
$arr[34,11] = 6000
$arr[34,12] = 6000
$arr[34,13] = "`$/docena de matas"
$arr[34,14] = "Región Metropolitana"
$arr[34,15] = 1000
$arr[34,16] = 6
$arr[34,17] = "Hortaliza"
$arr[35,0] = 4
$arr[35,1] = "Feria Lagunitas de Puerto Montt"
$arr[35,2] = "Los Lagos"
$arr[35,3] = 44201
$arr[35,4] = 10
$arr[35,5] = 100112052
$arr[35,6] = "Albahaca"
$arr[35,7] = "Sin especificar"
$arr[35,8] = "Primera"
$arr[35,9] = 70
$arr[35,10] = 7000
$arr[35,11] = 7000
$arr[35,12] = 7000
$arr[35,13] = "`$/docena de matas"
$arr[35,14] = "Región Metropolitana"
$arr[35,15] = 1167
$arr[35,16] = 6
$arr[35,17] = "Hortaliza"
$arr[36,0] = 4
$arr[36,1] = "Feria Lagunitas de Puerto Montt"
$arr[36,2] = "Los Lagos"
$arr[36,3] = 44579
$arr[36,4] = 10
$arr[36,5] = 100112052
$arr[36,6] = "Albahaca"
$arr[36,7] = "Sin especificar"
$arr[36,8] = "Primera"
$arr[36,9] = 120
$arr[36,10] = 7000
$arr[36,11] = 7000
$arr[36,12] = 7000
$arr[36,13] = "`$/docena de matas"
$arr[36,14] = "Región Metropolitana"
$arr[36,15] = 1167
$arr[36,16] = 6
$arr[36,17] = "Hortaliza"
$arr[37,0] = 4
$arr[37,1] = "Feria Lagunitas de Puerto Montt"
$arr[37,2] = "Los Lagos"
$arr[37,3] = 44225
$arr[37,4] = 10
$arr[37,5] = 100112052
$arr[37,6] = "Albahaca"
$arr[37,7] = "Sin especificar"
$arr[37,8] = "Primera"
$arr[37,9] = 120
$arr[37,10] = 6000
$arr[37,11] = 6000
$arr[37,12] = 6000
$arr[37,13] = "`$/docena de matas"
$arr[37,14] = "Región Metropolitana"
$arr[37,15] = 1000
$arr[37,16] = 6
$arr[37,17] = "Hortaliza"
$arr[38,0] = 4
$arr[38,1] = "Feria Lagunitas de Puerto Montt"
$arr[38,2] = "Los Lagos"
$arr[38,3] = 44252
$arr[38,4] = 10
$arr[38,5] = 100112052
$arr[38,6] = "Albahaca"
$arr[38,7] = "Sin especificar"
$arr[38,8] = "Primera"
$arr[38,9] = 60
$arr[38,10] = 6000
$arr[38,11] = 6000
$arr[38,12] = 6000
$arr[38,13] = "`$/docena de matas"
$arr[38,14] = "Región Metropolitana"
$arr[38,15] = 1000
$arr[38,16] = 6
$arr[38,17] = "Hortaliza"
$arr[39,0] = 4
$arr[39,1] = "Feria Lagunitas de Puerto Montt"
$arr[39,2] = "Los Lagos"
$arr[39,3] = 44271
$arr[39,4] = 10
$arr[39,5] = 100112052
$arr[39,6] = "Albahaca"
$arr[39,7] = "Sin especificar"
$arr[39,8] = "Primera"
$arr[39,9] = 90
$arr[39,10] = 7000
$arr[39,11] = 7000
$arr[39,12] = 7000
$arr[39,13] = "`$/docena de matas"
$arr[39,14] = "Región Metropolitana"
$arr[39,15] = 1167
$arr[39,16] = 6
$arr[39,17] = "Hortaliza"
$arr[40,0] = 4
$arr[40,1] = "Feria Lagunitas de Puerto Montt"
$arr[40,2] = "Los Lagos"
$arr[40,3] = 44162
$arr[40,4] = 10
$arr[40,5] = 100112052
$arr[40,6] = "Albahaca"
$arr[40,7] = "Sin especificar"
$arr[40,8] = "Primera"
$arr[40,9] = 90
$arr[40,10] = 6000
$arr[40,11] = 6000
$arr[40,12] = 6000
$arr[40,13] = "`$/docena de matas"
$arr[40,14] = "Región Metropolitana"
$arr[40,15] = 1000
$arr[40,16] = 6
$arr[40,17] = "Hortaliza"
$arr[41,0] = 4
$arr[41,1] = "Feria Lagunitas de Puerto Montt"
$arr[41,2] = "Los Lagos"
$arr[41,3] = 44516
$arr[41,4] = 10
$arr[41,5] = 100112052
$arr[41,6] = "Albahaca"
$arr[41,7] = "Sin especificar"
$arr[41,8] = "Primera"
$arr[41,9] = 100
$arr[41,10] = 8000
$arr[41,11] = 8000
$arr[41,12] = 8000
$arr[41,13] = "`$/docena de matas"
$arr[41,14] = "Región Metropolitana"
$arr[41,15] = 1333
$arr[41,16] = 6
$arr[41,17] = "Hortaliza"
$arr[42,0] = 4
$arr[42,1] = "Feria Lagunitas de Puerto Montt"
$arr[42,2] = "Los Lagos"
$arr[42,3] = 44568
$arr[42,4] = 10
$arr[42,5] = 100112052
$arr[42,6] = "Albahaca"
$arr[42,7] = "Sin especificar"
$arr[42,8] = "Primera"
$arr[42,9] = 80
$arr[42,10] = 9000
$arr[42,11] = 9000
$arr[42,12] = 9000
$arr[42,13] = "`$/docena de matas"
$arr[42,14] = "Región Metropolitana"
$arr[42,15] = 1500
$arr[42,16] = 6
$arr[42,17] = "Hortaliza"
$arr[43,0] = 4
$arr[43,1] = "Feria Lagunitas de Puerto Montt"
$arr[43,2] = "Los Lagos"
$arr[43,3] = 44231
$arr[43,4] = 10
$arr[43,5] = 100112052
$arr[43,6] = "Albahaca"
$arr[43,7] = "Sin especificar"
$arr[43,8] = "Primera"
$arr[43,9] = 40
$arr[43,10] = 6000
$arr[43,11] = 6000
$arr[43,12] = 6000
$arr[43,13] = "`$/docena de matas"
$arr[43,14] = "Región Metropolitana"
$arr[43,15] = 1000
$arr[43,16] = 6
$arr[43,17] = "Hortaliza"
$arr[44,0] = 4
$arr[44,1] = "Feria Lagunitas de Puerto Montt"
$arr[44,2] = "Los Lagos"
$arr[44,3] = 44565
$arr[44,4] = 10
$arr[44,5] = 100112052
$arr[44,6] = "Albahaca"
$arr[44,7] = "Sin especificar"
$arr[44,8] = "Primera"
$arr[44,9] = 80
$arr[44,10] = 7000
$arr[44,11] = 7000
$arr[44,12] = 7000
$arr[44,13] = "`$/docena de matas"
$arr[44,14] = "Región Metropolitana"
$arr[44,15] = 1167
$arr[44,16] = 6
$arr[44,17] = "Hortaliza"
$arr[45,0] = 4
$arr[45,1] = "Feria Lagunitas de Puerto Montt"
$arr[45,2] = "Los Lagos"
$arr[45,3] = 44334
$arr[45,4] = 10
$arr[45,5] = 100112052
$arr[45,6] = "Albahaca"
$arr[45,7] = "Sin especificar"
$arr[45,8] = "Primera"
$arr[45,9] = 120
$arr[45,10] = 6000
$arr[45,11] = 6000
$arr[45,12] = 6000
$arr[45,13] = "`$/docena de matas"
$arr[45,14] = "Región Metropolitana"
$arr[45,15] = 1000
$arr[45,16] = 6
$arr[45,17] = "Hortaliza"
$arr[46,0] = 4
$arr[46,1] = "Feria Lagunitas de Puerto Montt"
$arr[46,2] = "Los Lagos"
$arr[46,3] = 44280
$arr[46,4] = 10
$arr[46,5] = 100112052
$arr[46,6] = "Albahaca"
$arr[46,7] = "Sin especificar"
$arr[46,8] = "Primera"
$arr[46,9] = 80
$arr[46,10] = 6000
$arr[46,11] = 6000
$arr[46,12] = 6000
$arr[46,13] = "`$/docena de matas"
$arr[46,14] = "Región Metropolitana"
$arr[46,15] = 1000
$arr[46,16] = 6
$arr[46,17] = "Hortaliza"
$arr[47,0] = 4
$arr[47,1] = "Feria Lagunitas de Puerto Montt"
$arr[47,2] = "Los Lagos"
$arr[47,3] = 44567
$arr[47,4] = 10
$arr[47,5] = 100112052
$arr[47,6] = "Albahaca"
$arr[47,7] = "Sin especificar"
$arr[47,8] = "Primera"
$arr[47,9] = 80
$arr[47,10] = 7000
$arr[47,11] = 7000
$arr[47,12] = 7000
$arr[47,13] = "`$/docena de matas"
$arr[47,14] = "Región Metropolitana"
$arr[47,15] = 1167
$arr[47,16] = 6
$arr[47,17] = "Hortaliza"
$arr[48,0] = 4
$arr[48,1] = "Feria Lagunitas de Puerto Montt"
$arr[48,2] = "Los Lagos"
$arr[48,3] = 44642
$arr[48,4] = 10
$arr[48,5] = 100112052
$arr[48,6] = "Albahaca"
$arr[48,7] = "Sin especificar"
$arr[48,8] = "Primera"
$arr[48,9] = 90
$arr[48,10] = 5500
$arr[48,11] = 5500
$arr[48,12] = 5500
$arr[48,13] = "`$/docena de matas"
$arr[48,14] = "Región Metropolitana"
$arr[48,15] = 917
$arr[48,16] = 6
$arr[48,17] = "Hortaliza"
$arr[49,0] = 4
$arr[49,1] = "Feria Lagunitas de Puerto Montt"
$arr[49,2] = "Los Lagos"
$arr[49,3] = 44537
$arr[49,4] = 10
$arr[49,5] = 100112052
$arr[49,6] = "Albahaca"
$arr[49,7] = "Sin especificar"
$arr[49,8] = "Primera"
$arr[49,9] = 120
$arr[49,10] = 8000
$arr[49,11] = 8000
$arr[49,12] = 8000
$arr[49,13] = "`$/docena de matas"
$arr[49,14] = "Región Metropolitana"
$arr[49,15] = 1333
$arr[49,16] = 6
$arr[49,17] = "Hortaliza"
$arr[50,0] = 4
$arr[50,1] = "Feria Lagunitas de Puerto Montt"
$arr[50,2] = "Los Lagos"
$arr[50,3] = 44446
$arr[50,4] = 10
$arr[50,5] = 100112052
$arr[50,6] = "Albahaca"
$arr[50,7] = "Sin especificar"
$arr[50,8] = "Primera"
$arr[50,9] = 90
$arr[50,10] = 6000
$arr[50,11] = 6000
$arr[50,12] = 6000
$arr[50,13] = "`$/paquete"
$arr[50,14] = "Región de Arica y Parinacota"
$arr[50,15] = 6000
$arr[50,16] = 1
$arr[50,17] = "Hortaliza"
$arr[51,0] = 4
$arr[51,1] = "Feria Lagunitas de Puerto Montt"
$arr[51,2] = "Los Lagos"
$arr[51,3] = 44637
$arr[51,4] = 10
$arr[51,5] = 100112052
$arr[51,6] = "Albahaca"
$arr[51,7] = "Sin especificar"
$arr[51,8] = "Primera"
$arr[51,9] = 90
$arr[51,10] = 7000
$arr[51,11] = 7000
$arr[51,12] = 7000
$arr[51,13] = "`$/docena de matas"
$arr[51,14] = "Región Metropolitana"
$arr[51,15] = 1167
$arr[51,16] = 6
$arr[51,17] = "Hortaliza"
$arr[52,0] = 4
$arr[52,1] = "Feria Lagunitas de Puerto Montt"
$arr[52,2] = "Los Lagos"
$arr[52,3] = 44208
$arr[52,4] = 10
$arr[52,5] = 100112052
$arr[52,6] = "Albahaca"
$arr[52,7] = "Sin especificar"
$arr[52,8] = "Primera"
$arr[52,9] = 80
$arr[52,10] = 6500
$arr[52,11] = 6500
$arr[52,12] = 6500
$arr[52,13] = "`$/docena de matas"
$arr[52,14] = "Región Metropolitana"
$arr[52,15] = 1083
$arr[52,16] = 6
$arr[52,17] = "Hortaliza"
$arr[53,0] = 4
$arr[53,1] = "Feria Lagunitas de Puerto Montt"
$arr[53,2] = "Los Lagos"
$arr[53,3] = 44644
$arr[53,4] = 10
$arr[53,5] = 100112052
$arr[53,6] = "Albahaca"
$arr[53,7] = "Sin especificar"
$arr[53,8] = "Primera"
$arr[53,9] = 60
$arr[53,10] = 7000
$arr[53,11] = 7000
$arr[53,12] = 7000
$arr[53,13] = "`$/docena de matas"
$arr[53,14] = "Región Metropolitana"
$arr[53,15] = 1167
$arr[53,16] = 6
$arr[53,17] = "Hortaliza"
$arr[54,0] = 4
$arr[54,1] = "Feria Lagunitas de Puerto Montt"
$arr[54,2] = "Los Lagos"
$arr[54,3] = 44530
$arr[54,4] = 10
$arr[54,5] = 100112052
$arr[54,6] = "Albahaca"
$arr[54,7] = "Sin especificar"
$arr[54,8] = "Primera"
$arr[54,9] = 110
$arr[54,10] = 8000
$arr[54,11] = 8000
$arr[54,12] = 8000
$arr[54,13] = "`$/docena de matas"
$arr[54,14] = "Región Metropolitana"
$arr[54,15] = 1333
$arr[54,16] = 6
$arr[54,17] = "Hortaliza"
$arr[55,0] = 4
$arr[55,1] = "Feria Lagunitas de Puerto Montt"
$arr[55,2] = "Los Lagos"
$arr[55,3] = 44294
$arr[55,4] = 10
$arr[55,5] = 100112052
$arr[55,6] = "Albahaca"
$arr[55,7] = "Sin especificar"
$arr[55,8] = "Primera"
$arr[55,9] = 30
$arr[55,10] = 7000
$arr[55,11] = 7000
$arr[55,12] = 7000
$arr[55,13] = "`$/docena de matas"
$arr[55,14] = "Región Metropolitana"
$arr[55,15] = 1167
$arr[55,16] = 6
$arr[55,17] = "Hortaliza"
$arr[56,0] = 4
$arr[56,1] = "Feria Lagunitas de Puerto Montt"
$arr[56,2] = "Los Lagos"
$arr[56,3] = 44617
$arr[56,4] = 10
$arr[56,5] = 100112052
$arr[56,6] = "Albahaca"
$arr[56,7] = "Sin especificar"
$arr[56,8] = "Primera"
$arr[56,9] = 120
$arr[56,10] = 6000
$arr[56,11] = 6000
$arr[56,12] = 6000
$arr[56,13] = "`$/docena de matas"
$arr[56,14] = "Región Metropolitana"
$arr[56,15] = 1000
$arr[56,16] = 6
$arr[56,17] = "Hortaliza"
$arr[57,0] = 4
$arr[57,1] = "Feria Lagunitas de Puerto Montt"
$arr[57,2] = "Los Lagos"
$arr[57,3] = 44264
$arr[57,4] = 10
$arr[57,5] = 100112052
$arr[57,6] = "Albahaca"
$arr[57,7] = "Sin especificar"
$arr[57,8] = "Primera"
$arr[57,9] = 90
$arr[57,10] = 7000
$arr[57,11] = 7000
$arr[57,12] = 7000
$arr[57,13] = "`$/docena de matas"
$arr[57,14] = "Región Metropolitana"
$arr[57,15] = 1167
$arr[57,16] = 6
$arr[57,17] = "Hortaliza"
$arr[58,0] = 4
$arr[58,1] = "Feria Lagunitas de Puerto Montt"
$arr[58,2] = "Los Lagos"
$arr[58,3] = 44232
$arr[58,4] = 10
$arr[58,5] = 100112052
$arr[58,6] = "Albahaca"
$arr[58,7] = "Sin especificar"
$arr[58,8] = "Primera"
$arr[58,9] = 120
$arr[58,10] = 6000
$arr[58,11] = 6000
$arr[58,12] = 6000
$arr[58,13] = "`$/docena de matas"
$arr[58,14] = "Región Metropolitana"
$arr[58,15] = 1000
$arr[58,16] = 6
$arr[58,17] = "Hortaliza"
$arr[59,0] = 4
$arr[59,1] = "Feria Lagunitas de Puerto Montt"
$arr[59,2] = "Los Lagos"
$arr[59,3] = 44330
$arr[59,4] = 10
$arr[59,5] = 100112052
$arr[59,6] = "Albahaca"
$arr[59,7] = "Sin especificar"
$arr[59,8] = "Primera"
$arr[59,9] = 90
$arr[59,10] = 7000
$arr[59,11] = 7000
$arr[59,12] = 7000
$arr[59,13] = "`$/docena de matas"
$arr[59,14] = "Región Metropolitana"
$arr[59,15] = 1167
$arr[59,16] = 6
$arr[59,17] = "Hortaliza"
$arr[60,0] = 4
$arr[60,1] = "Feria Lagunitas de Puerto Montt"
$arr[60,2] = "Los Lagos"
$arr[60,3] = 44504
$arr[60,4] = 10
$arr[60,5] = 100112052
$arr[60,6] = "Albahaca"
$arr[60,7] = "Sin especificar"
$arr[60,8] = "Primera"
$arr[60,9] = 60
$arr[60,10] = 7000
$arr[60,11] = 7000
$arr[60,12] = 7000
$arr[60,13] = "`$/paquete"
$arr[60,14] = "Región de Arica y Parinacota"
$arr[60,15] = 7000
$arr[60,16] = 1
$arr[60,17] = "Hortaliza"
$arr[61,0] = 4
$arr[61,1] = "Feria Lagunitas de Puerto Montt"
$arr[61,2] = "Los Lagos"
$arr[61,3] = 44572
$arr[61,4] = 10
$arr[61,5] = 100112052
$arr[61,6] = "Albahaca"
$arr[61,7] = "Sin especificar"
$arr[61,8] = "Primera"
$arr[61,9] = 120
$arr[61,10] = 7000
$arr[61,11] = 7000
$arr[61,12] = 7000
$arr[61,13] = "`$/docena de matas"
$arr[61,14] = "Región Metropolitana"
$arr[61,15] = 1167
$arr[61,16] = 6
$arr[61,17] = "Hortaliza"
$arr[62,0] = 4
$arr[62,1] = "Feria Lagunitas de Puerto Montt"
$arr[62,2] = "Los Lagos"
$arr[62,3] = 44257
$arr[62,4] = 10
$arr[62,5] = 100112052
$arr[62,6] = "Albahaca"
$arr[62,7] = "Sin especificar"
$arr[62,8] = "Primera"
$arr[62,9] = 120
$arr[62,10] = 7000
$arr[62,11] = 7000
$arr[62,12] = 7000
$arr[62,13] = "`$/docena de matas"
$arr[62,14] = "Región Metropolitana"
$arr[62,15] = 1167
$arr[62,16] = 6
$arr[62,17] = "Hortaliza"
$arr[63,0] = 4
$arr[63,1] = "Feria Lagunitas de Puerto Montt"
$arr[63,2] = "Los Lagos"
$arr[63,3] = 44301
$arr[63,4] = 10
$arr[63,5] = 100112052
$arr[63,6] = "Albahaca"
$arr[63,7] = "Sin especificar"
$arr[63,8] = "Primera"
$arr[63,9] = 30
$arr[63,10] = 8000
$arr[63,11] = 8000
$arr[63,12] = 8000
$arr[63,13] = "`$/docena de matas"
$arr[63,14] = "Región Metropolitana"
$arr[63,15] = 1333
$arr[63,16] = 6
$arr[63,17] = "Hortaliza"
$arr[64,0] = 4
$arr[64,1] = "Feria Lagunitas de Puerto Montt"
$arr[64,2] = "Los Lagos"
$arr[64,3] = 44236
$arr[64,4] = 10
$arr[64,5] = 100112052
$arr[64,6] = "Albahaca"
$arr[64,7] = "Sin especificar"
$arr[64,8] = "Primera"
$arr[64,9] = 120
$arr[64,10] = 6000
$arr[64,11] = 6000
$arr[64,12] = 6000
$arr[64,13] = "`$/docena de matas"
$arr[64,14] = "Región Metropolitana"
$arr[64,15] = 1000
$arr[64,16] = 6
$arr[64,17] = "Hortaliza"
$arr[65,0] = 4
$arr[65,1] = "Feria Lagunitas de Puerto Montt"
$arr[65,2] = "Los Lagos"
$arr[65,3] = 44229
$arr[65,4] = 10
$arr[65,5] = 100112052
$arr[65,6] = "Albahaca"
$arr[65,7] = "Sin especificar"
$arr[65,8] = "Primera"
$arr[65,9] = 150
$arr[65,10] = 6000
$arr[65,11] = 6000
$arr[65,12] = 6000
$arr[65,13] = "`$/docena de matas"
$arr[65,14] = "Región Metropolitana"
$arr[65,15] = 1000
$arr[65,16] = 6
$arr[65,17] = "Hortaliza"
$arr[66,0] = 4
$arr[66,1] = "Feria Lagunitas de Puerto Montt"
$arr[66,2] = "Los Lagos"
$arr[66,3] = 44299
$arr[66,4] = 10
$arr[66,5] = 100112052
$arr[66,6] = "Albahaca"
$arr[66,7] = "Sin especificar"
$arr[66,8] = "Primera"
$arr[66,9] = 170
$arr[66,10] = 7000
$arr[66,11] = 8000
$arr[66,12] = 7471
$arr[66,13] = "`$/docena de matas"
$arr[66,14] = "Región Metropolitana"
$arr[66,15] = 1245
$arr[66,16] = 6
$arr[66,17] = "Hortaliza"
$arr[67,0] = 4
$arr[67,1] = "Feria Lagunitas de Puerto Montt"
$arr[67,2] = "Los Lagos"
$arr[67,3] = 44610
$arr[67,4] = 10
$arr[67,5] = 100112052
$arr[67,6] = "Albahaca"
$arr[67,7] = "Sin especificar"
$arr[67,8] = "Primera"
$arr[67,9] = 150
$arr[67,10] = 6000
$arr[67,11] = 6000
$arr[67,12] = 6000
$arr[67,13] = "`$/docena de matas"
$arr[67,14] = "Región Metropolitana"
$arr[67,15] = 1000
$arr[67,16] = 6
$arr[67,17] = "Hortaliza"
$arr[68,0] = 4
$arr[68,1] = "Feria Lagunitas de Puerto Montt"
$arr[68,2] = "Los Lagos"
$arr[68,3] = 44390
$arr[68,4] = 10
$arr[68,5] = 100112052
$arr[68,6] = "Albahaca"
$arr[68,7] = "Sin especificar"
$arr[68,8] = "Primera"
$arr[68,9] = 80
$arr[68,10] = 7000
$arr[68,11] = 7000
$arr[68,12] = 7000
$arr[68,13] = "`$/paquete"
$arr[68,14] = "Región de Arica y Parinacota"
$arr[68,15] = 7000
$arr[68,16] = 1
$arr[68,17] = "Hortaliza"
$arr[69,0] = 4
$arr[69,1] = "Feria Lagunitas de Puerto Montt"
$arr[69,2] = "Los Lagos"
$arr[69,3] = 44285
$arr[69,4] = 10
$arr[69,5] = 100112052
$arr[69,6] = "Albahaca"
$arr[69,7] = "Sin especificar"
$arr[69,8] = "Primera"
$arr[69,9] = 120
$arr[69,10] = 6000
$arr[69,11] = 6000
$arr[69,12] = 6000
$arr[69,13] = "`$/docena de matas"
$arr[69,14] = "Región Metropolitana"
$arr[69,15] = 1000
$arr[69,16] = 6
$arr[69,17] = "Hortaliza"
$arr[70,0] = 4
$arr[70,1] = "Feria Lagunitas de Puerto Montt"
$arr[70,2] = "Los Lagos"
$arr[70,3] = 44498
$arr[70,4] = 10
$arr[70,5] = 100112052
$arr[70,6] = "Albahaca"
$arr[70,7] = "Sin especificar"
$arr[70,8] = "Primera"
$arr[70,9] = 90
$arr[70,10] = 7000
$arr[70,11] = 7000
$arr[70,12] = 7000
$arr[70,13] = "`$/paquete"
$arr[70,14] = "Región de Arica y Parinacota"
$arr[70,15] = 7000
$arr[70,16] = 1
$arr[70,17] = "Hortaliza"
$arr[71,0] = 4
$arr[71,1] = "Feria Lagunitas de Puerto Montt"
$arr[71,2] = "Los Lagos"
$arr[71,3] = 44628
$arr[71,4] = 10
$arr[71,5] = 100112052
$arr[71,6] = "Albahaca"
$arr[71,7] = "Sin especificar"
$arr[71,8] = "Primera"
$arr[71,9] = 180
$arr[71,10] = 5000
$arr[71,11] = 6000
$arr[71,12] = 5500
$arr[71,13] = "`$/docena de matas"
$arr[71,14] = "Región Metropolitana"
$arr[71,15] = 917
$arr[71,16] = 6
$arr[71,17] = "Hortaliza"
$arr[72,0] = 4
$arr[72,1] = "Feria Lagunitas de Puerto Montt"
$arr[72,2] = "Los Lagos"
$arr[72,3] = 44544
$arr[72,4] = 10
$arr[72,5] = 100112052
$arr[72,6] = "Albahaca"
$arr[72,7] = "Sin especificar"
$arr[72,8] = "Primera"
$arr[72,9] = 80
$arr[72,10] = 7000
$arr[72,11] = 7000
$arr[72,12] = 7000
$arr[72,13] = "`$/docena de matas"
$arr[72,14] = "Región Metropolitana"
$arr[72,15] = 1167
$arr[72,16] = 6
$arr[72,17] = "Hortaliza"
$ws.Range("A39:R111").Value = $arr
$ws.Range("D39:D111").NumberFormat = "YYYY-MM-DD HH:MM:SS"
Write-Output "done"